# Atualizado por script em 26-11-2023 20:30
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data between row 36 and row 37 (columns F:V only; A:E stay the same) ---
$ws.Range("F36").Value = "Maghreb Fez"
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = "Jeunesse Sportive Soualem"
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 2.09
$ws.Range("K36").Value = "30/09/2023 06:42"
$ws.Range("L36").Value = 2.01
$ws.Range("M36").Value = "01/10/2023 19:13"
$ws.Range("N36").Value = 2.82
$ws.Range("O36").Value = "30/09/2023 06:42"
$ws.Range("P36").Value = 3.03
$ws.Range("Q36").Value = "01/10/2023 19:01"
$ws.Range("R36").Value = 3.56
$ws.Range("S36").Value = "30/09/2023 06:42"
$ws.Range("T36").Value = 4.13
$ws.Range("U36").Value = "01/10/2023 19:13"
$ws.Range("V36").Value = "https://www.betexplorer.com/football/morocco/botola-pro/maghreb-fez-jeunesse-sportive-soualem/OYJ1YDS0/"

$ws.Range("F37").Value = "Youssoufia Berrechid"
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = "Mouloudia Oujda"
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 2.28
$ws.Range("K37").Value = "01/10/2023 04:42"
$ws.Range("L37").Value = 2.05
$ws.Range("M37").Value = "01/10/2023 19:11"
$ws.Range("N37").Value = 2.9
$ws.Range("O37").Value = "01/10/2023 04:42"
$ws.Range("P37").Value = 2.99
$ws.Range("Q37").Value = "01/10/2023 19:11"
$ws.Range("R37").Value = 3.29
$ws.Range("S37").Value = "01/10/2023 04:42"
$ws.Range("T37").Value = 4.02
$ws.Range("U37").Value = "01/10/2023 19:11"
$ws.Range("V37").Value = "https://www.betexplorer.com/football/morocco/botola-pro/youssoufia-berrechid-mouloudia-oujda/IeJ5XXs7/"

# --- Append 5 new match rows (71-75), copying the formatting from row 70 first ---
$ws.Range("A70:V70").Copy()
$ws.Range("A71:V75").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 71
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "morocco"
$ws.Range("C71").Value = "botola-pro"
$ws.Range("D71").Value = "2023-2024"
$ws.Range("E71").Value = 45255.66666666666
$ws.Range("F71").Value = "Mouloudia Oujda"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Hassania Agadir"
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2.57
$ws.Range("K71").Value = "23/11/2023 15:12"
$ws.Range("L71").Value = 2.77
$ws.Range("M71").Value = "25/11/2023 15:55"
$ws.Range("N71").Value = 2.79
$ws.Range("O71").Value = "23/11/2023 15:12"
$ws.Range("P71").Value = 2.65
$ws.Range("Q71").Value = "25/11/2023 15:55"
$ws.Range("R71").Value = 2.73
$ws.Range("S71").Value = "23/11/2023 15:12"
$ws.Range("T71").Value = 3
$ws.Range("U71").Value = "25/11/2023 15:55"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/morocco/botola-pro/mouloudia-oujda-hassania-agadir/rLmIPU6Q/"

# Row 72
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = "morocco"
$ws.Range("C72").Value = "botola-pro"
$ws.Range("D72").Value = "2023-2024"
$ws.Range("E72").Value = 45255.76041666666
$ws.Range("F72").Value = "Youssoufia Berrechid"
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = "Maghreb Fez"
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2.89
$ws.Range("K72").Value = "23/11/2023 15:12"
$ws.Range("L72").Value = 3.62
$ws.Range("M72").Value = "25/11/2023 18:11"
$ws.Range("N72").Value = 2.82
$ws.Range("O72").Value = "23/11/2023 15:12"
$ws.Range("P72").Value = 2.87
$ws.Range("Q72").Value = "25/11/2023 18:11"
$ws.Range("R72").Value = 2.47
$ws.Range("S72").Value = "23/11/2023 15:12"
$ws.Range("T72").Value = 2.24
$ws.Range("U72").Value = "25/11/2023 18:11"
$ws.Range("V72").Value = "https://www.betexplorer.com/football/morocco/botola-pro/youssoufia-berrechid-maghreb-fez/pUZqLlFs/"

# Row 73
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = "morocco"
$ws.Range("C73").Value = "botola-pro"
$ws.Range("D73").Value = "2023-2024"
$ws.Range("E73").Value = 45256.625
$ws.Range("F73").Value = "Raja Casablanca"
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = "Chabab Mohammedia"
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1.33
$ws.Range("K73").Value = "25/11/2023 03:12"
$ws.Range("L73").Value = 1.39
$ws.Range("M73").Value = "26/11/2023 14:57"
$ws.Range("N73").Value = 4.36
$ws.Range("O73").Value = "25/11/2023 03:12"
$ws.Range("P73").Value = 4.1
$ws.Range("Q73").Value = "26/11/2023 14:57"
$ws.Range("R73").Value = 7.41
$ws.Range("S73").Value = "25/11/2023 03:12"
$ws.Range("T73").Value = 9.57
$ws.Range("U73").Value = "26/11/2023 14:57"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/morocco/botola-pro/raja-casablanca-chabab-mohammedia/zJSdIna0/"

# Row 74
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = "morocco"
$ws.Range("C74").Value = "botola-pro"
$ws.Range("D74").Value = "2023-2024"
$ws.Range("E74").Value = 45256.66666666666
$ws.Range("F74").Value = "IR Tanger"
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Moghreb Tetouan"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 2.06
$ws.Range("K74").Value = "25/11/2023 04:13"
$ws.Range("L74").Value = 2.43
$ws.Range("M74").Value = "26/11/2023 15:41"
$ws.Range("N74").Value = 2.96
$ws.Range("O74").Value = "25/11/2023 04:13"
$ws.Range("P74").Value = 2.9
$ws.Range("Q74").Value = "26/11/2023 15:44"
$ws.Range("R74").Value = 3.43
$ws.Range("S74").Value = "25/11/2023 04:13"
$ws.Range("T74").Value = 3.18
$ws.Range("U74").Value = "26/11/2023 15:41"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/morocco/botola-pro/ir-tanger-moghreb-tetouan/25UlK8Um/"

# Row 75
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "morocco"
$ws.Range("C75").Value = "botola-pro"
$ws.Range("D75").Value = "2023-2024"
$ws.Range("E75").Value = 45256.66666666666
$ws.Range("F75").Value = "Renaissance Zemamra"
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = "Jeunesse Sportive Soualem"
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1.98
$ws.Range("K75").Value = "25/11/2023 04:13"
$ws.Range("L75").Value = 2.38
$ws.Range("M75").Value = "26/11/2023 15:58"
$ws.Range("N75").Value = 3.02
$ws.Range("O75").Value = "25/11/2023 04:13"
$ws.Range("P75").Value = 2.71
$ws.Range("Q75").Value = "26/11/2023 15:58"
$ws.Range("R75").Value = 3.72
$ws.Range("S75").Value = "25/11/2023 04:13"
$ws.Range("T75").Value = 3.55
$ws.Range("U75").Value = "26/11/2023 15:58"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/morocco/botola-pro/renaissance-zemamra-jeunesse-sportive-soualem/tdThJSqf/"
